$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 327.63635
$ws.Range("I6").Value = 160.4
$ws.Range("K6").Value = 481.2
$ws.Range("M6").Value = -369.2

$ws.Range("H64").Value = 4616.5
$ws.Range("I64").Value = 5174.5
$ws.Range("J64").Value = 3500.5
$ws.Range("K64").Value = 5174.5
$ws.Range("L64").Value = 3500.5
$ws.Range("M64").Value = -4926.5
$ws.Range("N64").Value = -3996.5

$ws.Range("H67").Value = 4616.5
$ws.Range("I67").Value = 5174.5
$ws.Range("J67").Value = 3500.5
$ws.Range("K67").Value = 5174.5
$ws.Range("L67").Value = 3500.5
$ws.Range("M67").Value = -4316.5
$ws.Range("N67").Value = -5216.5

$ws.Range("H100").Value = 2654.2727
$ws.Range("I100").Value = 2469.7
$ws.Range("K100").Value = 2469.7
$ws.Range("M100").Value = -1928.7

$ws.Range("H107").Value = 571.26666
$ws.Range("I107").Value = 573.2143
$ws.Range("K107").Value = 573.2143
$ws.Range("M107").Value = 1346.7857

$ws.Range("H125").Value = 2160.6365
$ws.Range("I125").Value = 2049.5
$ws.Range("J125").Value = 2224.1428
$ws.Range("K125").Value = 18445.5
$ws.Range("L125").Value = 20017.2852
$ws.Range("M125").Value = -15985.5
$ws.Range("N125").Value = -24937.2852

$ws.Range("H129").Value = 2743.3333
$ws.Range("J129").Value = 4998.6665
$ws.Range("L129").Value = 14995.9995
$ws.Range("N129").Value = -24995.9995

$ws.Range("H132").Value = 62506570
$ws.Range("I132").Value = 62506570
$ws.Range("K132").Value = 187519710
$ws.Range("M132").Value = -187517180

$ws.Range("H135").Value = 1022.6429
$ws.Range("I135").Value = 1022.6429
$ws.Range("K135").Value = 9203.786100000001
$ws.Range("M135").Value = -6668.786100000001

$ws.Range("H137").Value = 25642928
$ws.Range("I137").Value = 37038710
$ws.Range("J137").Value = 2413.5
$ws.Range("K137").Value = 111116130
$ws.Range("L137").Value = 7240.5
$ws.Range("M137").Value = -111113580
$ws.Range("N137").Value = -12340.5

$ws.Range("H138").Value = 7124.4683
$ws.Range("I138").Value = 6777.273
$ws.Range("J138").Value = 7230.5557
$ws.Range("K138").Value = 20331.819
$ws.Range("L138").Value = 21691.6671
$ws.Range("M138").Value = -15191.819
$ws.Range("N138").Value = -31971.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15223
$ws.Range("I32").Value = 10586.134
$ws.Range("J32").Value = 49999.5
$ws.Range("K32").Value = 10586.134
$ws.Range("L32").Value = 49999.5
$ws.Range("M32").Value = -10299.134
$ws.Range("N32").Value = -50573.5

$ws.Range("H33").Value = 6666
$ws.Range("I33").Value = 6666
$ws.Range("K33").Value = 6666
$ws.Range("M33").Value = -6337

$ws.Range("H36").Value = 1506.5
$ws.Range("I36").Value = 1506.5
$ws.Range("K36").Value = 1506.5
$ws.Range("M36").Value = -1160.5

$ws.Range("H74").Value = 71269800
$ws.Range("I74").Value = 97184460
$ws.Range("K74").Value = 97184460
$ws.Range("M74").Value = -97183586

$ws.Range("H77").Value = 71269800
$ws.Range("I77").Value = 97184460
$ws.Range("K77").Value = 485922300
$ws.Range("M77").Value = -485917932

$ws.Range("H132").Value = 6645.684
$ws.Range("I132").Value = 6385.5293
$ws.Range("J132").Value = 8857
$ws.Range("K132").Value = 19156.5879
$ws.Range("L132").Value = 26571
$ws.Range("M132").Value = -16626.5879
$ws.Range("N132").Value = -31631

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 999
$ws.Range("I94").Value = 998.5
$ws.Range("J94").Value = 999.6667
$ws.Range("K94").Value = 998.5
$ws.Range("L94").Value = 999.6667
$ws.Range("M94").Value = -547.5
$ws.Range("N94").Value = -1901.6667

$ws.Range("H99").Value = 2484.1538
$ws.Range("I99").Value = 2510.7778
$ws.Range("J99").Value = 2424.25
$ws.Range("K99").Value = 2510.7778
$ws.Range("L99").Value = 2424.25
$ws.Range("M99").Value = -1012.7778
$ws.Range("N99").Value = -5420.25

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2693.0625
$ws.Range("I31").Value = 2693.0625
$ws.Range("K31").Value = 2693.0625
$ws.Range("M31").Value = -2398.0625

$ws.Range("H34").Value = 2693.0625
$ws.Range("I34").Value = 2693.0625
$ws.Range("K34").Value = 2693.0625
$ws.Range("M34").Value = -2491.0625

$ws.Range("H86").Value = 5249.625
$ws.Range("I86").Value = 5249.625
$ws.Range("K86").Value = 5249.625
$ws.Range("M86").Value = -4126.625

$ws.Range("H89").Value = 5249.625
$ws.Range("I89").Value = 5249.625
$ws.Range("K89").Value = 26248.125
$ws.Range("M89").Value = -20632.125

$ws.Range("H105").Value = 1867
$ws.Range("I105").Value = 1127.8334
$ws.Range("K105").Value = 1127.8334
$ws.Range("M105").Value = 619.1666

$ws.Range("H106").Value = 30671
$ws.Range("J106").Value = 30671
$ws.Range("L106").Value = 30671
$ws.Range("N106").Value = -33195

$ws.Range("H132").Value = 6672034.5
$ws.Range("I132").Value = 7148462
$ws.Range("K132").Value = 21445386
$ws.Range("M132").Value = -21442856

$ws.Range("H134").Value = 2697.4211
$ws.Range("I134").Value = 2349.5833
$ws.Range("J134").Value = 3293.7144
$ws.Range("K134").Value = 7048.749899999999
$ws.Range("L134").Value = 9881.143199999999
$ws.Range("M134").Value = -4513.749899999999
$ws.Range("N134").Value = -14951.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 507.44446
$ws.Range("I15").Value = 489
$ws.Range("K15").Value = 1467
$ws.Range("M15").Value = -1327

$ws.Range("H33").Value = 238.41667
$ws.Range("J33").Value = 229.25
$ws.Range("L33").Value = 1375.5
$ws.Range("N33").Value = -1941.5

$ws.Range("H44").Value = 143211.28
$ws.Range("I44").Value = 250044.75
$ws.Range("J44").Value = 766.6667
$ws.Range("K44").Value = 750134.25
$ws.Range("L44").Value = 2300.0001
$ws.Range("M44").Value = -749736.25
$ws.Range("N44").Value = -3096.0001

$ws.Range("H137").Value = 30474.75
$ws.Range("J137").Value = 9999.5
$ws.Range("L137").Value = 29998.5
$ws.Range("N137").Value = -40198.5

$ws.Range("H138").Value = 14521.667
$ws.Range("I138").Value = 14846
$ws.Range("J138").Value = 12900
$ws.Range("K138").Value = 44538
$ws.Range("L138").Value = 38700
$ws.Range("M138").Value = -39398
$ws.Range("N138").Value = -48980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H70").Value = 7998.5
$ws.Range("I70").Value = 7998.5
$ws.Range("K70").Value = 7998.5
$ws.Range("M70").Value = -7728.5

$ws.Range("H73").Value = 7998.5
$ws.Range("I73").Value = 7998.5
$ws.Range("K73").Value = 7998.5
$ws.Range("M73").Value = -7062.5

$ws.Range("H107").Value = 2018.25
$ws.Range("I107").Value = 878.1429000000001
$ws.Range("K107").Value = 878.1429000000001
$ws.Range("M107").Value = 1041.8571

$ws.Range("H113").Value = 10000
$ws.Range("J113").Value = 10000
$ws.Range("L113").Value = 10000
$ws.Range("N113").Value = -14340

$ws.Range("H132").Value = 17547696
$ws.Range("I132").Value = 4095.875
$ws.Range("J132").Value = 111113560
$ws.Range("K132").Value = 12287.625
$ws.Range("L132").Value = 333340680
$ws.Range("M132").Value = -9757.625
$ws.Range("N132").Value = -333345740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2001
$ws.Range("I2").Value = 2001
$ws.Range("K2").Value = 2001
$ws.Range("M2").Value = -1889

$ws.Range("H46").Value = 2170.6
$ws.Range("I46").Value = 200
$ws.Range("J46").Value = 2663.25
$ws.Range("K46").Value = 200
$ws.Range("L46").Value = 2663.25
$ws.Range("M46").Value = -12
$ws.Range("N46").Value = -3039.25

$ws.Range("H132").Value = 2545
$ws.Range("J132").Value = 2200
$ws.Range("L132").Value = 6600
$ws.Range("N132").Value = -11660

$ws.Range("H136").Value = 2483
$ws.Range("J136").Value = 3449
$ws.Range("L136").Value = 10347
$ws.Range("N136").Value = -15447

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376

$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H122").Value = 3072.5789
$ws.Range("I122").Value = 3180
$ws.Range("K122").Value = 9540
$ws.Range("M122").Value = -7090

$ws.Range("H132").Value = 142858060
$ws.Range("I132").Value = 1075.3334
$ws.Range("K132").Value = 3226.0002
$ws.Range("M132").Value = -696.0001999999999

$ws.Range("H136").Value = 6906.56
$ws.Range("I136").Value = 7219.136
$ws.Range("J136").Value = 4614.3335
$ws.Range("K136").Value = 21657.408
$ws.Range("L136").Value = 13843.0005
$ws.Range("M136").Value = -19107.408
$ws.Range("N136").Value = -18943.0005
